# Updates cryptocurrency price/volume/hour data (symbol list refresh).
# Cell text must remain literal strings (e.g. "244.54", "-0.57%", "13"),
# matching the source file's inlineStr cell style -- so we force the
# NumberFormat to Text ("@") before assigning, which stops Excel from
# re-interpreting numeric-looking / percent-looking text as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.54'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.57%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '13'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '27.08'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '3.44%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '13'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.148'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.43%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '13'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05655'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.00%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '13'
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.473'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-0.20%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '13'
# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.85%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '13'
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8396'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.19%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '13'
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1333'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-0.68%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '13'
# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.52%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '13'
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02958'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '5.80%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '13'
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09401'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.02%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '13'
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001512'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.98%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '13'
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04283'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-8.54%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '13'
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0005970'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.65%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '13'
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006192'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.38%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '13'
# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-1.36%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '13'
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.006'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-1.15%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '13'
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.225'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '5.04%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '13'
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3113'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '13'
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.03142'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.93%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '13'
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1256'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.45%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '13'
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.591'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.89%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '13'
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1374'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.07%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '13'
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001222'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.02%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '13'
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004469'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.95%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '13'
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009800'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.12%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '13'
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.00007257'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-47.77%'
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '13'
# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '13'
# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '13'
# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '13'
# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '13'
# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '13'
# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '13'
# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '13'
# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '13'
# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '13'
# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '13'
# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '13'
# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.24%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '13'
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006050'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.90%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '13'
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1055'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.17%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '13'
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002610'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '0.63%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '13'
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008371'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-4.43%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '13'
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005364'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.39%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '13'
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.03%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '13'
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1100'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-17.16%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '13'
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '28.65%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '13'
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.03%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '13'
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.03%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '13'
# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '13'

Write-Output "Updated 119 cells across 50 rows."
